$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.033.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.413.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '560.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.09%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.511'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.167'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.410.60'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.332'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '68.929.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000176'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.855.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '23.80'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.403.91'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '341.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.85'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.79'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.537.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0842'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.33'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.75%  '
$ws.Range('E32').Value = '  +11.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '449.03'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.55%  '
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.46%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  +5.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.12'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.301'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.91%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.77'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.37'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.51%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.51'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.48%  '
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '134.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.37'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0721'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.485'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.46%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0930'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.38%  '
